$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Relocate the hidden "_GoBack" bookmark so it sits right after the
#    "适合去看樱花" paragraph (this mirrors where Word leaves the mark
#    after the edit below is made). A directly-collapsed Range cannot
#    be used with Bookmarks.Add reliably, so we temporarily insert a
#    one-character placeholder, wrap the bookmark around it, and then
#    delete the placeholder text - leaving a collapsed bookmark exactly
#    where the placeholder used to be.
# ------------------------------------------------------------------
$anchorParaText = $d.Content
$anchorParaText.Find.Execute("适合去看樱花", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$placeholderStart = $anchorParaText.End
$placeholder = $d.Range($placeholderStart, $placeholderStart)
$placeholder.InsertAfter("X")

$wrap = $d.Range($placeholderStart, $placeholderStart + 1)
$d.Bookmarks.Add("_GoBack", $wrap)

$shrink = $d.Range($placeholderStart, $placeholderStart + 1)
$shrink.Text = ""

# ------------------------------------------------------------------
# 2) Extend the last diary entry with the new sentence. The existing
#    text is split across two runs (with the old _GoBack bookmark in
#    between); searching across the merged text and replacing it
#    collapses it back into a single run carrying the full sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "下雨，今天学习了分支管理，创建了一个dev分支",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "下雨，今天学习了分支管理，创建了一个dev分支。使用Git创建分支简单又快速。",
    2
)
